# Journal_de_travail.xlsx - record two new work sessions (45 + 40 min) on
# 2024-05-23 for the "Implémentation" subject, and start a new (still open)
# session right after. Matches commit:
#   "Implémentation de la génération du mouvement des pièces restantes"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 76: close out the existing entry (was Début-only, Fin missing) ---
# Début moves from 08:00 to 08:50, and a Fin of 09:35 is now recorded.
$ws.Range("C76").Value = 0.36805555555555558
$ws.Range("D76").Value = 0.39930555555555558
$ws.Range("D75").Copy() | Out-Null
$ws.Range("D76").PasteSpecial(-4122) | Out-Null
$ws.Range("G76").Value = "Résolution d'erreur de la génération des mouvements des pièces"

# --- Row 77: brand-new entry, 09:50 -> 10:30 ---
$ws.Range("B77").Value = 45435
$ws.Range("C77").Value = 0.40972222222222227
$ws.Range("D77").Value = 0.4375
$ws.Range("B76").Copy() | Out-Null
$ws.Range("B77").PasteSpecial(-4122) | Out-Null
$ws.Range("C76").Copy() | Out-Null
$ws.Range("C77").PasteSpecial(-4122) | Out-Null
$ws.Range("D76").Copy() | Out-Null
$ws.Range("D77").PasteSpecial(-4122) | Out-Null
$ws.Range("F77").Value = "Implémentation"
$ws.Range("G77").Value = "Implémentation de la génération du mouvement des pièces restante"
$ws.Rows.Item(77).RowHeight = 30

# --- Row 78: brand-new, still-open entry starting at 10:30 (no Fin yet) ---
$ws.Range("B78").Value = 45435
$ws.Range("C78").Value = 0.4375
$ws.Range("B76").Copy() | Out-Null
$ws.Range("B78").PasteSpecial(-4122) | Out-Null
$ws.Range("C76").Copy() | Out-Null
$ws.Range("C78").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Update the view so it reflects where work is now happening ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G78:G79").Select() | Out-Null
